# Update cryptocurrency price/volume data per Fri Jul 14 2023 run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.778.66"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.958.54"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  +3.13%  "
$ws.Range("D5").Value = "'250.46"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "'0.7041"
$ws.Range("E6").Value = "  -19.55%  "
$ws.Range("E7").Value = "  +2.64%  "
$ws.Range("D8").Value = "'0.3293"
$ws.Range("E8").Value = "  -3.91%  "
$ws.Range("D9").Value = "'26.97"
$ws.Range("E9").Value = "  +4.38%  "
$ws.Range("D10").Value = "'0.06925"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "'0.8094"
$ws.Range("E11").Value = "  -5.86%  "
$ws.Range("D12").Value = "'0.08028"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "1.963.23"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").Value = "'5.465"
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").Value = "'95.68"
$ws.Range("E15").Value = "  -5.82%  "
$ws.Range("D16").Value = "'14.86"
$ws.Range("E16").Value = "  +7.05%  "
$ws.Range("D17").Value = "'267.15"
$ws.Range("E17").Value = "  -4.55%  "
$ws.Range("D18").Value = "30.796.89"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").Value = "'5.970"
$ws.Range("E19").Value = "  +5.29%  "
$ws.Range("D20").Value = "'0.000007957"
$ws.Range("E20").Value = "  +1.32%  "
$ws.Range("D21").Value = "2.236.63"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").Value = "'1.005"
$ws.Range("E22").Value = "  +2.13%  "
$ws.Range("D23").Value = "'1.008"
$ws.Range("E23").Value = "  +2.91%  "
$ws.Range("D24").Value = "'6.941"
$ws.Range("E24").Value = "  +2.98%  "
$ws.Range("D25").Value = "'9.785"
$ws.Range("E25").Value = "  +2.30%  "
$ws.Range("D26").Value = "'160.58"
$ws.Range("D27").Value = "'19.24"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("D28").Value = "'2.307"
$ws.Range("E28").Value = "  +4.82%  "
$ws.Range("D29").Value = "'0.1301"
$ws.Range("E29").Value = "  -29.39%  "
$ws.Range("D30").Value = "'1.378"
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("D31").Value = "'1.574"
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("D32").Value = "'4.510"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").Value = "'4.295"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").Value = "'0.05191"
$ws.Range("E34").Value = "  +1.88%  "
$ws.Range("D35").Value = "'1.240"
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("D36").Value = "'0.7632"
$ws.Range("E36").Value = "  +3.23%  "
$ws.Range("D37").Value = "'2.771"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("D38").Value = "'0.01962"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").Value = "'2.859"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").Value = "'81.59"
$ws.Range("E40").Value = "  +4.98%  "
$ws.Range("D41").Value = "'6.609"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").Value = "'0.4518"
$ws.Range("E42").Value = "  -2.65%  "
$ws.Range("D43").Value = "'2.067"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").Value = "'0.8463"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("D45").Value = "'1.007"
$ws.Range("E45").Value = "  +2.44%  "
$ws.Range("D46").Value = "'102.60"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("D47").Value = "'9.978"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").Value = "'7.435"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").Value = "'36.58"
$ws.Range("E49").Value = "  +1.17%  "
$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").Value = "'2.896"
$ws.Range("E50").Value = "  +35.16%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.523"
$ws.Range("E51").Value = "  +6.68%  "
